# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Sun Feb 11 08:48:56 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "48.191.04"
$ws.Cells.Item(2, 5).Value = "  +2.17%  "
$ws.Cells.Item(3, 4).Value = "2.529.68"
$ws.Cells.Item(3, 5).Value = "  +1.51%  "
$ws.Cells.Item(4, 5).Value = "  -0.02%  "
$ws.Cells.Item(5, 4).Value = "'324.17"
$ws.Cells.Item(5, 5).Value = "  +0.49%  "
$ws.Cells.Item(6, 4).Value = "'109.20"
$ws.Cells.Item(6, 5).Value = "  +0.45%  "
$ws.Cells.Item(7, 5).Value = "  +0.46%  "
$ws.Cells.Item(8, 5).Value = "  +0.00%  "
$ws.Cells.Item(9, 5).Value = "  +4.20%  "
$ws.Cells.Item(10, 4).Value = "'40.92"
$ws.Cells.Item(10, 5).Value = "  +4.90%  "
$ws.Cells.Item(11, 4).Value = "'20.48"
$ws.Cells.Item(11, 5).Value = "  +11.27%  "
$ws.Cells.Item(12, 4).Value = "'0.0828"
$ws.Cells.Item(13, 5).Value = "  +1.35%  "
$ws.Cells.Item(14, 4).Value = "'7.31"
$ws.Cells.Item(14, 5).Value = "  +1.55%  "
$ws.Cells.Item(15, 4).Value = "2.927.78"
$ws.Cells.Item(15, 5).Value = "  +1.49%  "
$ws.Cells.Item(16, 4).Value = "2.532.37"
$ws.Cells.Item(16, 5).Value = "  +1.45%  "
$ws.Cells.Item(17, 4).Value = "'0.860"
$ws.Cells.Item(17, 5).Value = "  +0.93%  "
$ws.Cells.Item(18, 4).Value = "48.042.29"
$ws.Cells.Item(18, 5).Value = "  +1.99%  "
$ws.Cells.Item(19, 4).Value = "'13.31"
$ws.Cells.Item(19, 5).Value = "  +4.53%  "
$ws.Cells.Item(20, 4).Value = "'6.65"
$ws.Cells.Item(20, 5).Value = "  +0.58%  "
$ws.Cells.Item(21, 5).Value = "  +1.59%  "
$ws.Cells.Item(22, 5).Value = "  -0.27%  "
$ws.Cells.Item(23, 4).Value = "'72.30"
$ws.Cells.Item(23, 5).Value = "  +2.31%  "
$ws.Cells.Item(24, 4).Value = "'270.28"
$ws.Cells.Item(24, 5).Value = "  +9.23%  "
$ws.Cells.Item(25, 5).Value = "  -0.31%  "
$ws.Cells.Item(26, 4).Value = "'26.23"
$ws.Cells.Item(26, 5).Value = "  +1.17%  "
$ws.Cells.Item(27, 5).Value = "  -0.27%  "
$ws.Cells.Item(28, 2).Value = "Toncoin"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(28, 4).Value = "'2.29"
$ws.Cells.Item(28, 5).Value = "  +0.06%  "
$ws.Cells.Item(29, 2).Value = "Cosmos"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(29, 4).Value = "'10.16"
$ws.Cells.Item(29, 5).Value = "  +0.89%  "
$ws.Cells.Item(30, 5).Value = "  +0.81%  "
$ws.Cells.Item(31, 4).Value = "'35.70"
$ws.Cells.Item(31, 5).Value = "  +1.23%  "
$ws.Cells.Item(32, 4).Value = "'49.72"
$ws.Cells.Item(32, 5).Value = "  -0.39%  "
$ws.Cells.Item(33, 5).Value = "  -0.63%  "
$ws.Cells.Item(34, 5).Value = "  +0.21%  "
$ws.Cells.Item(35, 5).Value = "  -0.06%  "
$ws.Cells.Item(36, 5).Value = "  +1.27%  "
$ws.Cells.Item(37, 5).Value = "  +1.18%  "
$ws.Cells.Item(38, 4).Value = "'4.75"
$ws.Cells.Item(38, 5).Value = "  +1.37%  "
$ws.Cells.Item(39, 5).Value = "  +1.07%  "
$ws.Cells.Item(40, 4).Value = "'0.112"
$ws.Cells.Item(40, 5).Value = "  +0.13%  "
$ws.Cells.Item(41, 4).Value = "'22.25"
$ws.Cells.Item(41, 5).Value = "  +4.61%  "
$ws.Cells.Item(42, 4).Value = "'119.65"
$ws.Cells.Item(42, 5).Value = "  -1.38%  "
$ws.Cells.Item(43, 5).Value = "  -1.21%  "
$ws.Cells.Item(44, 5).Value = "  +1.72%  "
$ws.Cells.Item(45, 4).Value = "2.015.59"
$ws.Cells.Item(46, 4).Value = "'3.16"
$ws.Cells.Item(46, 5).Value = "  +3.97%  "
$ws.Cells.Item(47, 5).Value = "  +6.30%  "
$ws.Cells.Item(48, 5).Value = "  -0.10%  "
$ws.Cells.Item(49, 5).Value = "  +0.48%  "
$ws.Cells.Item(50, 5).Value = "  +1.69%  "
$ws.Cells.Item(51, 4).Value = "'79.73"
$ws.Cells.Item(51, 5).Value = "  +2.41%  "
